# Applies the price/volume refresh captured in the commit diff
# ("Updated cryptos list ... with GitHub Actions") to Sheet1 of the
# cryptos workbook. Columns D (Price) and E (Volume(1h)) are plain
# text cells in the source OOXML (t="inlineStr"), so any value that
# Excel would otherwise auto-detect as a number (e.g. "1.002") is
# force-formatted as Text first to keep it a literal string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.165.02'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '1.825.64'
$ws.Range('E3').Value = '  +0.83%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.63'
$ws.Range('E5').Value = '  -0.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('E7').Value = '  -3.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3927'
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09866'
$ws.Range('E9').Value = '  +25.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.109'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.21'
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.441'
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.000'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = '1.825.69'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.302'
$ws.Range('E16').Value = '  -0.81%  '
$ws.Range('E17').Value = '  +5.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.65'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06645'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.23'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.998'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').Value = '28.225.31'
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.33'
$ws.Range('E24').Value = '  +1.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.240'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.84'
$ws.Range('E26').Value = '  -1.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.73'
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').Value = '2.033.48'
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.417'
$ws.Range('E29').Value = '  +0.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.11'
$ws.Range('E30').Value = '  -1.30%  '
$ws.Range('E31').Value = '  -3.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.041'
$ws.Range('E32').Value = '  -2.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.599'
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.611'
$ws.Range('E34').Value = '  -1.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.06752'
$ws.Range('E35').Value = '  -6.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.036'
$ws.Range('E36').Value = '  -1.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02343'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2157'
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.981'
$ws.Range('E39').Value = '  -1.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.41'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6232'
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.179'
$ws.Range('E42').Value = '  +1.57%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.15'
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5940'
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.710'
$ws.Range('E46').Value = '  -0.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.272'
$ws.Range('E47').Value = '  -2.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.09'
$ws.Range('E48').Value = '  -1.36%  '
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.182'
$ws.Range('E50').Value = '  -3.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06792'
$ws.Range('E51').Value = '  -0.88%  '
